$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (SCD0171 -> SCD0011)
$ws.Name = "SCD0011"

# Update the test case id text in B2 (DGS-186 -> SCD0011-002)
$ws.Range("B2").Value = "SCD0011-002"

# Reflect the author's cursor/selection ending up on B2 after the edit
$ws.Range("B2").Select() | Out-Null
